$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kardex Cliente")

# --- Text / label updates (shared strings) ---
$ws.Range("B1").Value = "Estuche Spigen Samsung Galaxy Note 8 - Negro"
$ws.Range("E2").Value = "2C"
$ws.Range("B3").Value = "587CS22051"

# --- Header numeric updates ---
$ws.Range("B2").Value = 8808522199474
$ws.Range("E3").Value = 5

# --- Movement rows: update row 5, rewrite row 6, delete old rows 7 & 8 ---
$ws.Range("A5").Value = 44158
$ws.Range("B5").Value = 0.014375
$ws.Range("C5").Value = 4172472869

$ws.Range("A6").Value = 44155
$ws.Range("B6").Value = 0.64954861111111
$ws.Range("C6").ClearContents()
$ws.Range("D6").Value = 6
$ws.Range("E6").ClearContents()

# Remove the two trailing data rows (7 and 8) entirely so the used range
# shrinks back down to A1:E6.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(7).Delete()

# Reset the active selection to match the new last row.
$ws.Range("A6").Select()
